$d = $word.ActiveDocument

# The paragraph currently has the leading text split across two runs:
#   "This is a text with" + " "
# and the trailing text split across two runs:
#   " " + "insertions."
# Re-typing the text spanning each run boundary (as Word's Find & Replace
# does when editing) merges each pair back into a single run, matching the
# golden layout:
#   "This is a text with " + <ins>two exciting</ins> + " insertions."

$findReplace1 = $d.Content
$findReplace1.Find.Execute("This is a text with ", $false, $false, $false, $false, $false, $true, 1, $false, "This is a text with ", 2)

$findReplace2 = $d.Content
$findReplace2.Find.Execute(" insertions.", $false, $false, $false, $false, $false, $true, 1, $false, " insertions.", 2)
